$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.528.56'
$ws.Range("E2").Value = '  +5.09%  '
$ws.Range("D3").Value = '1.602.09'
$ws.Range("E3").Value = '  +2.76%  '
$ws.Range("E4").Value = '  -0.27%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '214.99'
$ws.Range("E5").Value = '  +2.42%  '
$ws.Range("E6").Value = '  +1.80%  '
$ws.Range("E7").Value = '  -0.08%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '24.04'
$ws.Range("E8").Value = '  +9.26%  '
$ws.Range("E9").Value = '  +1.58%  '
$ws.Range("E10").Value = '  +1.18%  '
$ws.Range("E11").Value = '  +2.31%  '
$ws.Range("D12").Value = '1.831.10'
$ws.Range("E12").Value = '  +2.85%  '
$ws.Range("D13").Value = '1.605.20'
$ws.Range("E13").Value = '  +2.86%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.535'
$ws.Range("E14").Value = '  +3.62%  '
$ws.Range("E15").Value = '  +0.79%  '
$ws.Range("D16").Value = '28.529.98'
$ws.Range("E16").Value = '  +5.26%  '
$ws.Range("E17").Value = '  +2.72%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '232.15'
$ws.Range("E18").Value = '  +7.47%  '
$ws.Range("E19").Value = '  +1.54%  '
$ws.Range("E20").Value = '  +1.69%  '
$ws.Range("E21").Value = '  -0.17%  '
$ws.Range("E22").Value = '  +0.64%  '
$ws.Range("E23").Value = '  +2.55%  '
$ws.Range("E24").Value = '  +1.90%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '152.54'
$ws.Range("E25").Value = '  +0.06%  '
$ws.Range("E26").Value = '  +2.04%  '
$ws.Range("E27").Value = '  +0.39%  '
$ws.Range("E28").Value = '  +1.21%  '
$ws.Range("E29").Value = '  -0.12%  '
$ws.Range("E30").Value = '  +1.31%  '
$ws.Range("E31").Value = '  +1.41%  '
$ws.Range("E32").Value = '  +1.15%  '
$ws.Range("E33").Value = '  +0.72%  '
$ws.Range("D34").Value = '1.423.34'
$ws.Range("E34").Value = '  -0.78%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.61'
$ws.Range("E35").Value = '  -0.12%  '
$ws.Range("E36").Value = '  -4.16%  '
$ws.Range("E37").Value = '  -0.13%  '
$ws.Range("E38").Value = '  +1.44%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.545'
$ws.Range("E39").Value = '  +2.61%  '
$ws.Range("E40").Value = '  +8.16%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.823'
$ws.Range("E41").Value = '  +2.30%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.76'
$ws.Range("E42").Value = '  -2.38%  '
$ws.Range("E43").Value = '  -0.14%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.985'
$ws.Range("E44").Value = '  -1.12%  '
$ws.Range("E45").Value = '  +6.57%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '65.01'
$ws.Range("E46").Value = '  +1.50%  '
$ws.Range("D47").Value = '1.741.20'
$ws.Range("E47").Value = '  +2.82%  '
$ws.Range("E48").Value = '  +0.83%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '87.56'
$ws.Range("E49").Value = '  +2.61%  '
$ws.Range("E50").Value = '  +9.56%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0528'
$ws.Range("E51").Value = '  +0.79%  '
